$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns in order: B, C, D, E, F, I, J, K, L, M, N (G stays 1, H is empty, A is row index - unchanged)
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @(
    @(2,"1.02","1.03215926294198","1.038960079507079","1.031625385671663","1.044770361999965","1.02558616476744","1.037291017445208","1.041746626773963","1.034432962556984","1.047540477808525","1.016353115855461"),
    @(3,"1.02","1.034186319873802","1.040886284846098","1.033370803401584","1.046987738301942","1.025828105562731","1.038954780536555","1.043480010256753","1.035984474811181","1.049565497943743","1.016903889678523"),
    @(4,"1.02","1.035484549232085","1.042117532369615","1.034488971503151","1.048398825102347","1.025973924201731","1.040017877053239","1.044586124264488","1.036976763452705","1.050851833588921","1.017255752221955"),
    @(5,"1.02","1.036027170756988","1.042631582726102","1.034956407287077","1.048986449222006","1.026032670262067","1.040461628462656","1.045047478366854","1.037391179267411","1.051386940368189","1.017402608661475"),
    @(6,"1.02","1.036118095947067","1.042717686548636","1.035034738217075","1.049084787895521","1.026042384461634","1.040535951549264","1.045124728891773","1.037460601970349","1.051476456814669","1.01742720443144"),
    @(7,"1.02","1.035491812086495","1.042124415072708","1.034495227722536","1.048406698854255","1.025974719195748","1.040023818888607","1.044592303195455","1.036982311613451","1.050859005889565","1.017257718693914"),
    @(8,"1.02","1.032847128319707","1.039614219128583","1.032217612967016","1.045524690968382","1.025670160678914","1.037856113896958","1.042335675213977","1.034959742089841","1.04822985545198","1.016540199667627"),
    @(9,"1.02","1.028081467259077","1.035072292009155","1.028115897825299","1.040260978749667","1.025050615002052","1.033930808343035","1.03823789906243","1.031304410892354","1.043409719032736","1.015240395542701"),
    @(10,"1.02","1.024829588530007","1.031960576969148","1.025318805641453","1.036621815910532","1.0245808267817","1.03123949493674","1.035420670281251","1.028803006031989","1.040065194280118","1.014348870685105"),
    @(11,"1.02","1.023402830685783","1.030592354780063","1.024092009384632","1.035013874004046","1.024363681300849","1.030055642519903","1.034179621828111","1.027703826576688","1.038584637146462","1.013956629237239"),
    @(12,"1.02","1.022869980814924","1.030080921286713","1.023633903771137","1.034411662809632","1.024280937526924","1.029613053579488","1.033715378696254","1.027293063561802","1.038029718140071","1.01380997628149"),
    @(13,"1.02","1.022984411007672","1.030190772405078","1.023732279592517","1.034541065117313","1.024298781166054","1.029708120772116","1.033815109418631","1.027381286902692","1.038148977122556","1.013841477555743"),
    @(14,"1.02","1.023358844568125","1.030550145619896","1.024054191990297","1.034964196723039","1.024356884381359","1.030019116634314","1.034141314388498","1.027669923743079","1.038538869515384","1.013944526528493"),
    @(15,"1.02","1.023589160139907","1.030771138690418","1.02425221018314","1.035224242510542","1.024392406479638","1.030210351046604","1.03434186520611","1.027847431906104","1.038778432479117","1.01400789082733"),
    @(16,"1.02","1.024923877327264","1.032050935319169","1.025399888750709","1.03672784218053","1.024594946818987","1.031317667345023","1.035502581652197","1.028875611069936","1.040162762762528","1.014374769634492"),
    @(17,"1.02","1.025756053426938","1.032848082316095","1.026115562638548","1.037662318683578","1.024718304018293","1.032007254141838","1.036224944054085","1.0295162162452","1.041022376006485","1.014603224583008"),
    @(18,"1.02","1.026239651816549","1.033311040967002","1.026531499248702","1.038204288313743","1.024788933752946","1.032407699093065","1.036644248991897","1.029888326044569","1.041520659947653","1.014735881501053"),
    @(19,"1.02","1.026404244222369","1.033468560568427","1.026673069862599","1.038388564527374","1.024812793071723","1.032543940869497","1.036786878080189","1.030014945895228","1.041690037056171","1.014781013554133"),
    @(20,"1.02","1.025666955149128","1.032762763845646","1.026038933647884","1.037562379134471","1.024705205914672","1.031933452580785","1.036147652663612","1.029447645586123","1.040930470527291","1.014578775453625"),
    @(21,"1.02","1.023248663639237","1.030444408617871","1.023959464159756","1.034839732633946","1.024339832228271","1.029927615490076","1.034045345893486","1.027584996367237","1.0384241941714","1.013914207797245"),
    @(22,"1.02","1.021711431994304","1.028968121176213","1.022637987549708","1.033099194749899","1.024098023431437","1.02864991739968","1.032704622981153","1.026399497778746","1.036819559516304","1.013490817618975"),
    @(23,"1.02","1.022527964574351","1.029752526093623","1.023339881425849","1.034024648483641","1.024227364850584","1.029328843815372","1.033417187303246","1.027029338925668","1.03767297961879","1.013715799419071"),
    @(24,"1.02","1.025707220370197","1.032801321776564","1.026073563633362","1.03760754707066","1.024711128468898","1.031966805823475","1.036182583605139","1.029478634483688","1.040972008272782","1.014589824806088"),
    @(25,"1.02","1.029326385316588","1.036260947339985","1.029187076272365","1.041644271901987","1.025220704692519","1.034958448226621","1.039312028701415","1.032260534160157","1.044678546865346","1.015580741984272"),
)

foreach ($rowEntry in $data) {
    $r = $rowEntry[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $col = $cols[$i]
        $val = [double]$rowEntry[$i + 1]
        $ws.Cells.Item($r, $col).Value = $val
    }
}
